$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the current row 161 (shifts old 161..166 down to 163..168)
$ws.Rows("161:162").Insert()

# New row 161 data
$ws.Cells.Item(161, 1).Value = 3
$ws.Cells.Item(161, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(161, 3).Value = "Coquimbo"
$ws.Cells.Item(161, 4).Value = 44448
$ws.Cells.Item(161, 5).Value = 5
$ws.Cells.Item(161, 6).Value = 100112013
$ws.Cells.Item(161, 7).Value = "Alcachofa"
$ws.Cells.Item(161, 8).Value = "Argentina(o)"
$ws.Cells.Item(161, 9).Value = "Primera"
$ws.Cells.Item(161, 10).Value = 150
$ws.Cells.Item(161, 11).Value = 10500
$ws.Cells.Item(161, 12).Value = 11000
$ws.Cells.Item(161, 13).Value = 10733
$ws.Cells.Item(161, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(161, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(161, 16).Value = 215
$ws.Cells.Item(161, 17).Value = 50
$ws.Cells.Item(161, 18).Value = "Hortaliza"

# New row 162 data
$ws.Cells.Item(162, 1).Value = 3
$ws.Cells.Item(162, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(162, 3).Value = "Coquimbo"
$ws.Cells.Item(162, 4).Value = 44448
$ws.Cells.Item(162, 5).Value = 5
$ws.Cells.Item(162, 6).Value = 100112013
$ws.Cells.Item(162, 7).Value = "Alcachofa"
$ws.Cells.Item(162, 8).Value = "Española"
$ws.Cells.Item(162, 9).Value = "Extra"
$ws.Cells.Item(162, 10).Value = 170
$ws.Cells.Item(162, 11).Value = 11500
$ws.Cells.Item(162, 12).Value = 12000
$ws.Cells.Item(162, 13).Value = 11765
$ws.Cells.Item(162, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(162, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(162, 16).Value = 392
$ws.Cells.Item(162, 17).Value = 30
$ws.Cells.Item(162, 18).Value = "Hortaliza"

# Match the date format used by the rest of column D (numFmtId 165, "YYYY-MM-DD HH:MM:SS")
$ws.Cells.Item(161, 4).NumberFormat = $ws.Cells.Item(163, 4).NumberFormat
$ws.Cells.Item(162, 4).NumberFormat = $ws.Cells.Item(163, 4).NumberFormat
